# "Add files via upload" — adds a new teammate's initials/file reference
# to the "Features" sheet (row 6: "Resize the full bracket without side
# scrolling.") and leaves that sheet active/selected, matching the author's
# last on-screen state when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Features")

# Fill in the new "Who" / "Which files are touched" values for row 6.
$ws.Range("C6").Value = "Thien"
$ws.Range("D6").Value = "MarchMadnessGUI"

# Make "Features" the active sheet and leave the selection where the author
# left it (one cell to the right of the newly entered data).
$ws.Activate()
$ws.Range("E6").Select() | Out-Null
